$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "I2" = 1001.125
    "M2" = -888.125
    "K2" = 1001.125
    "H2" = 895
    "I31" = 215.85715
    "L31" = 240
    "M31" = -417.5714499999999
    "K31" = 647.5714499999999
    "H31" = 175.1
    "J31" = 80
    "N31" = -700
    "I38" = 216.625
    "N38" = -6757.5
    "L38" = 6013.5
    "M38" = -277.875
    "K38" = 649.875
    "H38" = 415.27777
    "J38" = 2004.5
    "I62" = 81092.53999999999
    "M62" = -80468.53999999999
    "K62" = 81092.53999999999
    "H62" = 40303.75
    "I65" = 81092.53999999999
    "M65" = -402342.7
    "K65" = 405462.7
    "H65" = 40303.75
    "I132" = 1446.4348
    "M132" = -1809.3044
    "K132" = 4339.3044
    "H132" = 1432.3846
    "I137" = 3080.389
    "N137" = -22068.1758
    "L137" = 16968.1758
    "M137" = -6691.167000000001
    "K137" = 9241.167000000001
    "H137" = 3906.547
    "J137" = 5656.0586
    "I138" = 1672.6842
    "N138" = -18308.45
    "L138" = 8028.450000000001
    "M138" = 121.9474
    "K138" = 5018.0526
    "H138" = 2353
    "J138" = 2676.15
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "I23" = 0
    "L23" = 0
    "K23" = 0
    "H23" = 0
    "J23" = 0
    "I32" = 8076.516
    "M32" = -7789.516
    "K32" = 8076.516
    "H32" = 36921.31
    "N63" = -3294
    "L63" = 1922
    "H63" = 1784.1818
    "J63" = 1922
    "N66" = -16474
    "L66" = 9610
    "H66" = 1784.1818
    "J66" = 1922
    "I74" = 1330.8148
    "M74" = -456.8148000000001
    "K74" = 1330.8148
    "H74" = 1468.8788
    "I77" = 1330.8148
    "M77" = -2286.074000000001
    "K77" = 6654.074000000001
    "H77" = 1468.8788
    "N80" = -21993.5
    "L80" = 19997.5
    "H80" = 19997.5
    "J80" = 19997.5
    "N83" = -69976.5
    "L83" = 59992.5
    "H83" = 19997.5
    "J83" = 19997.5
    "I132" = 1658.3334
    "M132" = -2445.0002
    "K132" = 4975.0002
    "H132" = 1706.5714
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
$clears = @("N23", "M23")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "N82" = -25764
    "L82" = 24998
    "H82" = 14518.846
    "J82" = 24998
    "N85" = -27650
    "L85" = 24998
    "H85" = 14518.846
    "J85" = 24998
    "I134" = 1411.3846
    "N134" = -14595
    "L134" = 9525
    "M134" = -1699.1538
    "K134" = 4234.1538
    "H134" = 1646.5333
    "J134" = 3175
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "I31" = 1966.1538
    "M31" = -1671.1538
    "K31" = 1966.1538
    "H31" = 2769.9062
    "I34" = 1966.1538
    "M34" = -1764.1538
    "K34" = 1966.1538
    "H34" = 2769.9062
    "N68" = -26495.5
    "L68" = 24997.5
    "H68" = 24664.445
    "J68" = 24997.5
    "N71" = -82480.5
    "L71" = 74992.5
    "H71" = 24664.445
    "J71" = 24997.5
    "N74" = -79697.5
    "L74" = 77949.5
    "H74" = 59373.668
    "J74" = 77949.5
    "N75" = -44065
    "L75" = 42069
    "H75" = 42069
    "J75" = 42069
    "N77" = -242584.5
    "L77" = 233848.5
    "H77" = 59373.668
    "J77" = 77949.5
    "N78" = -136191
    "L78" = 126207
    "H78" = 42069
    "J78" = 42069
    "I93" = 17400
    "N93" = -34117.666
    "L93" = 30373.666
    "M93" = -15528
    "K93" = 17400
    "H93" = 22265.125
    "J93" = 30373.666
    "I107" = 1776.2
    "M107" = 143.8
    "K107" = 1776.2
    "H107" = 1701.2106
    "I132" = 2256.442
    "M132" = -4239.326
    "K132" = 6769.326
    "H132" = 2291.2183
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "I69" = 5873.8335
    "N69" = -23064.5
    "L69" = 21442.5
    "M69" = -16810.5005
    "K69" = 17621.5005
    "H69" = 6192.25
    "J69" = 7147.5
    "I72" = 5873.8335
    "N72" = -72439.5
    "L72" = 64327.5
    "M72" = -48808.5015
    "K72" = 52864.5015
    "H72" = 6192.25
    "J72" = 7147.5
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "N36" = -30968
    "L36" = 29998
    "H36" = 29998
    "J36" = 29998
    "I102" = 1434.5883
    "N102" = -4077
    "L102" = 833
    "M102" = 187.4117000000001
    "K102" = 1434.5883
    "H102" = 1344.35
    "J102" = 833
    "I107" = 601.5
    "M107" = 1318.5
    "K107" = 601.5
    "H107" = 55558096
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "I22" = 756.2308
    "M22" = -461.2308
    "K22" = 756.2308
    "H22" = 764.6087
    "I27" = 756.2308
    "M27" = -649.2308
    "K27" = 756.2308
    "H27" = 764.6087
    "I46" = 39419.637
    "M46" = -39231.637
    "K46" = 39419.637
    "H46" = 25198.111
    "I68" = 1683.5
    "M68" = -934.5
    "K68" = 1683.5
    "H68" = 2260.375
    "I71" = 1683.5
    "M71" = -4673.5
    "K71" = 8417.5
    "H71" = 2260.375
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "N96" = -4112
    "L96" = 1366
    "H96" = 3525.875
    "J96" = 1366
    "I132" = 2203446
    "M132" = -6607808
    "K132" = 6610338
    "H132" = 3491632.8
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
